# The authored change swaps the content of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: the presentation's main (document) theme, which was
# "Integral" / "Red Violet" becomes the default "Office Theme" / "Office"
# palette (and vice-versa for the notes-master theme part).
#
# The PowerPoint object model only exposes write access to the single
# document theme (ppt/theme/theme1.xml, the theme used by the slide master /
# the actual slides) via Slide.ThemeColorScheme(...).Colors(i).RGB - there is
# no supported object-model path that reaches the notes-master's theme part,
# so we recolor the document theme to match the target "Office" palette.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function Set-ThemeColor($scheme, [int]$index, [byte]$r, [byte]$g, [byte]$b) {
    $val = ($b * 65536) + ($g * 256) + $r
    $scheme.Colors($index).RGB = $val
}

# MsoThemeColorSchemeIndex order: 1=dk1 2=lt1 3=dk2 4=lt2
# 5=accent1 6=accent2 7=accent3 8=accent4 9=accent5 10=accent6
# 11=hlink 12=folHlink
Set-ThemeColor $tcs 1  0x00 0x00 0x00   # dk1
Set-ThemeColor $tcs 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor $tcs 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor $tcs 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor $tcs 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor $tcs 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor $tcs 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor $tcs 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor $tcs 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor $tcs 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor $tcs 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor $tcs 12 0x95 0x4F 0x72   # folHlink
